$d = $word.ActiveDocument

# --- Change 1 ---------------------------------------------------------
# In the "ADICIONALES" row of the first table the quantity cell reads
# "4 TURNOS". The leading "4" is removed, leaving just the space that
# already separated it from "TURNOS" (so the cell becomes " TURNOS").
$d.Content.Find.Execute("4 TURNOS", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " TURNOS", 2)

# --- Change 2 ---------------------------------------------------------
# Between the "ANEXOS" table and the "NOMBRE / C.C" table there are two
# consecutive, completely empty paragraphs (just a paragraph mark with
# the Arial rPr, no text runs). One of them is removed so only a single
# blank paragraph separates the two tables.
$anexosTable = $d.Tables.Item(2)
$nombreTable = $d.Tables.Item(3)

$gapStart = $anexosTable.Range.End
$gapEnd = $nombreTable.Range.Start
$gap = $d.Range($gapStart, $gapEnd)

if (($gapEnd - $gapStart) -eq 2 -and $gap.Text -eq "`r`r") {
    # Delete just the first of the two paragraph marks.
    $d.Range($gapStart, $gapStart + 1).Delete()
}
